$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "Exeter" college, keeping the alphabetical ordering ---
# (it currently sits, out of place, as the very last row of the table)
$ws.Rows("5:5").Insert() | Out-Null
$ws.Rows("5:5").RowHeight = 19

$ws.Range("A5").Value = "Exeter"
$ws.Range("B5").Value = "Exeter College"
$ws.Range("C5").Value = "Turl St"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "Oxford "
$ws.Range("F5").Value = "OX1 3DP"

# The Exeter address cell picked up a slightly different (near-identical) dark-grey
# Arial font when it was retyped.
$ws.Range("C5").Font.Name = "Arial"
$ws.Range("C5").Font.Size = 14
$ws.Range("C5").Font.Color = 2039583

# --- Remove the old, out-of-place "Exeter" row further down the table ---
# (after the insert above it has shifted from row 20 down to row 21)
$ws.Rows("21:21").Delete() | Out-Null

# Deleting that row pulled everything below it back up by one row, but the
# stray trailing cell near the bottom of the sheet actually ends up one row
# lower than before overall, so compensate by opening up a blank row there.
$ws.Rows("28:28").Insert() | Out-Null

# --- Small text tweak: the "St Catherine's" address gained a trailing space ---
$ws.Range("C15").Value = "Manor Rd "

# --- Restore the active selection shown in the saved workbook ---
$ws.Range("A20").Select() | Out-Null
